$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the dataset. It belongs right
# above the existing row 515, so insert a new row there (this pushes the
# old rows 515-619 down to 516-620) and then fill in the new record's
# values.
$ws.Rows.Item(515).Insert()

$ws.Range("A515").Value = 3
$ws.Range("B515").Value = "Femacal de La Calera"
$ws.Range("C515").Value = "Coquimbo"
$ws.Range("D515").Value = 45209
$ws.Range("E515").Value = 5
$ws.Range("F515").Value = 100112012
$ws.Range("G515").Value = "Espinaca"
$ws.Range("H515").Value = "Sin especificar"
$ws.Range("I515").Value = "Primera"
$ws.Range("J515").Value = 110
$ws.Range("K515").Value = 4000
$ws.Range("L515").Value = 4000
$ws.Range("M515").Value = 4000
$ws.Range("N515").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O515").Value = "Provincia de Quillota"
$ws.Range("P515").Value = 1333
$ws.Range("Q515").Value = 3
$ws.Range("R515").Value = "Hortaliza"
